$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2990.6936
$ws.Range("I76").Value = 2661.8728
$ws.Range("J76").Value = 5574.2856
$ws.Range("K76").Value = 2661.8728
$ws.Range("L76").Value = 5574.2856
$ws.Range("M76").Value = -2346.8728
$ws.Range("N76").Value = -6204.2856

$ws.Range("H79").Value = 2990.6936
$ws.Range("I79").Value = 2661.8728
$ws.Range("J79").Value = 5574.2856
$ws.Range("K79").Value = 2661.8728
$ws.Range("L79").Value = 5574.2856
$ws.Range("M79").Value = -1569.8728
$ws.Range("N79").Value = -7758.2856

$ws.Range("H137").Value = 1538.0571
$ws.Range("I137").Value = 1530.9048
$ws.Range("J137").Value = 1548.7858
$ws.Range("K137").Value = 4592.7144
$ws.Range("L137").Value = 4646.357400000001
$ws.Range("M137").Value = -2042.7144
$ws.Range("N137").Value = -9746.357400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4658.1685
$ws.Range("I32").Value = 4077.4521
$ws.Range("K32").Value = 4077.4521
$ws.Range("M32").Value = -3790.4521

$ws.Range("H61").Value = 1566.6842
$ws.Range("I61").Value = 1369.1428
$ws.Range("J61").Value = 2119.8
$ws.Range("K61").Value = 1369.1428
$ws.Range("L61").Value = 2119.8
$ws.Range("M61").Value = -1157.1428
$ws.Range("N61").Value = -2543.8

$ws.Range("H122").Value = 1852.8485
$ws.Range("I122").Value = 1509.1111
$ws.Range("J122").Value = 3399.6667
$ws.Range("K122").Value = 4527.3333
$ws.Range("L122").Value = 10199.0001
$ws.Range("M122").Value = -2077.3333
$ws.Range("N122").Value = -15099.0001

$ws.Range("H136").Value = 1566.6842
$ws.Range("I136").Value = 1369.1428
$ws.Range("J136").Value = 2119.8
$ws.Range("K136").Value = 4107.428400000001
$ws.Range("L136").Value = 6359.400000000001
$ws.Range("M136").Value = -1557.428400000001
$ws.Range("N136").Value = -11459.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1615.5
$ws.Range("J134").Value = 2677.3333
$ws.Range("L134").Value = 8031.999899999999
$ws.Range("N134").Value = -13101.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 437.125
$ws.Range("I22").Value = 148.75
$ws.Range("J22").Value = 725.5
$ws.Range("K22").Value = 148.75
$ws.Range("L22").Value = 725.5
$ws.Range("M22").Value = 201.25
$ws.Range("N22").Value = -1425.5

$ws.Range("H31").Value = 4069.5908
$ws.Range("I31").Value = 6003
$ws.Range("K31").Value = 6003
$ws.Range("M31").Value = -5708

$ws.Range("H34").Value = 4069.5908
$ws.Range("I34").Value = 6003
$ws.Range("K34").Value = 6003
$ws.Range("M34").Value = -5801

$ws.Range("H122").Value = 3135.3333
$ws.Range("I122").Value = 1870.6666
$ws.Range("J122").Value = 4400
$ws.Range("K122").Value = 5611.9998
$ws.Range("L122").Value = 13200
$ws.Range("M122").Value = -3161.9998
$ws.Range("N122").Value = -18100

$ws.Range("H134").Value = 2444.9678
$ws.Range("I134").Value = 2230
$ws.Range("J134").Value = 3562.8
$ws.Range("K134").Value = 6690
$ws.Range("L134").Value = 10688.4
$ws.Range("M134").Value = -4155
$ws.Range("N134").Value = -15758.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 878474.5600000001
$ws.Range("I5").Value = 781
$ws.Range("J5").Value = 1951211.1
$ws.Range("K5").Value = 2343
$ws.Range("L5").Value = 5853633.300000001
$ws.Range("M5").Value = -2231
$ws.Range("N5").Value = -5853857.300000001

$ws.Range("H81").Value = 4507
$ws.Range("J81").Value = 4507
$ws.Range("L81").Value = 13521
$ws.Range("N81").Value = -15767

$ws.Range("H84").Value = 4507
$ws.Range("J84").Value = 4507
$ws.Range("L84").Value = 40563
$ws.Range("N84").Value = -51795

$ws.Range("H121").Value = 34255.566
$ws.Range("I121").Value = 333599.34
$ws.Range("J121").Value = 995.14813
$ws.Range("K121").Value = 1000798.02
$ws.Range("L121").Value = 2985.44439
$ws.Range("M121").Value = -999488.02
$ws.Range("N121").Value = -5605.444390000001

$ws.Range("H129").Value = 1830.138
$ws.Range("I129").Value = 793.125
$ws.Range("J129").Value = 2225.1904
$ws.Range("K129").Value = 2379.375
$ws.Range("L129").Value = 6675.5712
$ws.Range("M129").Value = 2620.625
$ws.Range("N129").Value = -16675.5712

$ws.Range("H130").Value = 1583.25
$ws.Range("I130").Value = 1499.5
$ws.Range("K130").Value = 4498.5
$ws.Range("M130").Value = 521.5

$ws.Range("H131").Value = 7469.4116
$ws.Range("I131").Value = 875
$ws.Range("K131").Value = 2625
$ws.Range("M131").Value = 2415

$ws.Range("H135").Value = 878474.5600000001
$ws.Range("I135").Value = 781
$ws.Range("J135").Value = 1951211.1
$ws.Range("K135").Value = 7029
$ws.Range("L135").Value = 17560899.9
$ws.Range("M135").Value = -4494
$ws.Range("N135").Value = -17565969.9

$ws.Range("H139").Value = 1558.6666
$ws.Range("I139").Value = 1098.8572
$ws.Range("J139").Value = 2053.8462
$ws.Range("K139").Value = 3296.5716
$ws.Range("L139").Value = 6161.5386
$ws.Range("M139").Value = 1843.4284
$ws.Range("N139").Value = -16441.5386

$ws.Range("H140").Value = 2153
$ws.Range("I140").Value = 1763.45
$ws.Range("J140").Value = 2752.3076
$ws.Range("K140").Value = 5290.35
$ws.Range("L140").Value = 8256.9228
$ws.Range("M140").Value = -110.3500000000004
$ws.Range("N140").Value = -18616.9228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 2330.3845
$ws.Range("I132").Value = 2112.6843
$ws.Range("J132").Value = 2921.2856
$ws.Range("K132").Value = 6338.0529
$ws.Range("L132").Value = 8763.856800000001
$ws.Range("M132").Value = -3808.0529
$ws.Range("N132").Value = -13823.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -705
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -893
$ws.Range("N27").ClearContents()

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 19356.5
$ws.Range("J108").Value = 19356.5
$ws.Range("L108").Value = 19356.5
$ws.Range("N108").Value = -27036.5

$ws.Range("H132").Value = 1079.1613
$ws.Range("I132").Value = 926.7349
$ws.Range("J132").Value = 2344.3
$ws.Range("K132").Value = 2780.2047
$ws.Range("L132").Value = 7032.900000000001
$ws.Range("M132").Value = -250.2047000000002
$ws.Range("N132").Value = -12092.9
